$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the numeric-looking "Price" text values are written as text,
# not auto-converted to numbers by Excel (matches source data which
# stores prices such as "70.706.80" or "1.00" as literal text).
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

# Apply the updated cell values from the refreshed cryptos feed.
$ws.Range('D2').Value = '70.706.80'
$ws.Range('E2').Value = '  +5.43%  '
$ws.Range('D3').Value = '3.651.72'
$ws.Range('E3').Value = '  +5.71%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '593.51'
$ws.Range('E5').Value = '  +1.27%  '
$ws.Range('D6').Value = '194.41'
$ws.Range('E6').Value = '  +4.08%  '
$ws.Range('D7').Value = '0.647'
$ws.Range('E7').Value = '  +2.68%  '
$ws.Range('D8').Value = '3.645.99'
$ws.Range('E8').Value = '  +5.77%  '
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('D10').Value = '0.180'
$ws.Range('E10').Value = '  +4.48%  '
$ws.Range('E11').Value = '  +4.66%  '
$ws.Range('D12').Value = '58.40'
$ws.Range('E12').Value = '  +3.86%  '
$ws.Range('D13').Value = '0.0000293'
$ws.Range('E13').Value = '  +5.83%  '
$ws.Range('D14').Value = '9.95'
$ws.Range('E14').Value = '  +5.95%  '
$ws.Range('D15').Value = '4.225.50'
$ws.Range('E15').Value = '  +5.48%  '
$ws.Range('E16').Value = '  +6.67%  '
$ws.Range('D17').Value = '3.642.16'
$ws.Range('E17').Value = '  +5.69%  '
$ws.Range('D18').Value = '70.646.64'
$ws.Range('E18').Value = '  +5.56%  '
$ws.Range('D19').Value = '12.80'
$ws.Range('E19').Value = '  +5.72%  '
$ws.Range('E20').Value = '  +3.21%  '
$ws.Range('E21').Value = '  +4.18%  '
$ws.Range('D22').Value = '490.14'
$ws.Range('E22').Value = '  +0.07%  '
$ws.Range('D23').Value = '19.14'
$ws.Range('E23').Value = '  +14.71%  '
$ws.Range('D24').Value = '5.32'
$ws.Range('E24').Value = '  -0.77%  '
$ws.Range('E25').Value = '  +1.79%  '
$ws.Range('D26').Value = '91.18'
$ws.Range('E26').Value = '  +1.70%  '
$ws.Range('D27').Value = '3.19'
$ws.Range('E27').Value = '  +8.63%  '
$ws.Range('D28').Value = '11.64'
$ws.Range('E28').Value = '  +6.51%  '
$ws.Range('D29').Value = '9.65'
$ws.Range('E29').Value = '  +6.72%  '
$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').Value = '7.95'
$ws.Range('E30').Value = '  +11.77%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '32.89'
$ws.Range('E31').Value = '  +5.17%  '
$ws.Range('D32').Value = '0.121'
$ws.Range('E32').Value = '  +8.93%  '
$ws.Range('D33').Value = '630.22'
$ws.Range('E33').Value = '  +5.65%  '
$ws.Range('D34').Value = '12.28'
$ws.Range('E34').Value = '  +4.84%  '
$ws.Range('D35').Value = '65.74'
$ws.Range('E35').Value = '  +2.93%  '
$ws.Range('D36').Value = '40.71'
$ws.Range('E36').Value = '  +11.16%  '
$ws.Range('B37').Value = 'TheGraph'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D37').Value = '0.414'
$ws.Range('E37').Value = '  +7.89%  '
$ws.Range('B38').Value = 'PEPE'
$ws.Range('C38').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D38').Value = '0.0₃0825'
$ws.Range('E38').Value = '  +9.55%  '
$ws.Range('D39').Value = '0.148'
$ws.Range('E39').Value = '  -2.11%  '
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('D41').Value = '3.59'
$ws.Range('E41').Value = '  +1.67%  '
$ws.Range('D42').Value = '3.305.89'
$ws.Range('E42').Value = '  +1.78%  '
$ws.Range('D43').Value = '2.87'
$ws.Range('E43').Value = '  +14.22%  '
$ws.Range('E44').Value = '  +9.33%  '
$ws.Range('D45').Value = '0.0454'
$ws.Range('E45').Value = '  +5.87%  '
$ws.Range('D46').Value = '2.93'
$ws.Range('E46').Value = '  +3.94%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = '3.32'
$ws.Range('E47').Value = '  +2.65%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').Value = '0.139'
$ws.Range('E48').Value = '  +3.08%  '
$ws.Range('B49').Value = 'THORChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D49').Value = '9.26'
$ws.Range('E49').Value = '  +6.40%  '
$ws.Range('E50').Value = '  +1.29%  '
$ws.Range('D51').Value = '0.998'
$ws.Range('E51').Value = '  +0.00%  '
